# Add "RMSLE.TissuePC" and "N.TissuePC" columns to the httk-benchmarks
# table, inserted right before the existing "Notes" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item(1)

# Insert two new blank worksheet columns at the position of the current
# "Notes" column (P). This shifts the "Notes" column (data + style) two
# places to the right, to column R.
$ws.Range("P1").EntireColumn.Insert()
$ws.Range("P1").EntireColumn.Insert()

# The new columns should keep the same (narrow) width as the other
# RMSLE/N columns (K:O) instead of the default width.
$ws.Range("P1:Q1").EntireColumn.ColumnWidth = $ws.Columns.Item(11).ColumnWidth

# Grow the table definition to cover the two new columns.
$tbl.Resize($ws.Range("A1:R23"))

# Set the new header names (writing directly to the header cells, which
# also updates the underlying table column metadata). Re-assert the
# "Notes" header too, since it now occupies a freshly-grown table slot.
$ws.Range("P1").Value = "RMSLE.TissuePC"
$ws.Range("Q1").Value = "N.TissuePC"
$ws.Range("R1").Value = "Notes"

# Populate the RMSLE.TissuePC / N.TissuePC data cells.
$ws.Range("P4").Value = 0.4612
$ws.Range("Q4").Value = 12

$ws.Range("P5").Value = 0.4612
$ws.Range("Q5").Value = 12

$ws.Range("P6").Value = 0.5563
$ws.Range("Q6").Value = 412

$ws.Range("P7").Value = 0.5925
$ws.Range("Q7").Value = 964

$ws.Range("P8").Value = 0.5926
$ws.Range("Q8").Value = 964

$ws.Range("P9").Value = 0.5925
$ws.Range("Q9").Value = 964

$ws.Range("P10").Value = 0.6136
$ws.Range("Q10").Value = 953

$ws.Range("P11").Value = 0.6136
$ws.Range("Q11").Value = 953

$ws.Range("P12").Value = 0.6136
$ws.Range("Q12").Value = 953

$ws.Range("P13").Value = 0.6115
$ws.Range("Q13").Value = 964

$ws.Range("P14").Value = 0.6115
$ws.Range("Q14").Value = 964

$ws.Range("P15").Value = 0.6098
$ws.Range("Q15").Value = 858

$ws.Range("P16").Value = 0.7611
$ws.Range("Q16").Value = 858

$ws.Range("P17").Value = 0.7611
$ws.Range("Q17").Value = 858

$ws.Range("P18").Value = 0.7854
$ws.Range("Q18").Value = 851

$ws.Range("P19").Value = 0.7866
$ws.Range("Q19").Value = 840

$ws.Range("P20").Value = 0.5995
$ws.Range("Q20").Value = 863

$ws.Range("P21").Value = 0.6428
$ws.Range("Q21").Value = 863

$ws.Range("P22").Value = 0.643
$ws.Range("Q22").Value = 863

$ws.Range("P23").Value = 0.63
$ws.Range("Q23").Value = 863

# Match the final selection state left behind in the source workbook
# (the view had also been scrolled so column D is leftmost).
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("R4").Select()
